# CYBER360-Ex-3.2-CIM.docx — "Add files via upload" edit
#
# Real, semantically meaningful changes made by the author (the rest of
# the raw XML diff is noise from Word re-saving the file with a newer
# build: proofing-tag churn from a fresh spell/grammar pass, a refreshed
# theme-font table, and new xml namespace declarations):
#
#   1. The cached SAVEDATE field result in the "Last Updated:" line was
#      refreshed to the new save timestamp.
#   2. The "Title" content control's lock changed from "sdtLocked"
#      (structure locked only) to "sdtContentLocked" (structure AND
#      contents locked) — matching the "Subject" control next to it.
#   3. The "Task 2" heading text changed from
#      "Task 2—Getting CIM-Instances" to "Task 2—Getting CIM Instances".
#   4. The "Task 3" heading text changed from "Task 3—WMI Queries" to
#      "Task 3—WMI Queries with CIM".

$d = $word.ActiveDocument
$emDash = [char]0x2014

# 1) Refresh the cached "Last Updated" SAVEDATE field text.
$d.Content.Find.Execute(
    "1/9/2024 3:49 PM", $true, $false, $false, $false, $false,
    $true, 1, $false, "4/18/2024 9:10 AM", 2) | Out-Null

# 2) Lock the "Title" content control's contents too (sdtLocked ->
#    sdtContentLocked), same as the existing "Subject" control.
#    Identify it precisely by its (stable) w:id / alias="Title".
$ccs = $d.ContentControls
for ($i = 1; $i -le $ccs.Count; $i++) {
    $cc = $ccs.Item($i)
    if ($cc.Title -eq "Title" -and $cc.ID -eq "1406417971") {
        $cc.LockContents = $true
    }
}

# 3) "Task 2—Getting CIM-Instances" -> "Task 2—Getting CIM Instances"
$d.Content.Find.Execute(
    ($emDash + "Getting CIM-Instances"), $true, $false, $false, $false,
    $false, $true, 1, $false, ($emDash + "Getting CIM Instances"), 2) | Out-Null

# 4) "Task 3—WMI Queries" -> "Task 3—WMI Queries with CIM"
$d.Content.Find.Execute(
    ("Task 3" + $emDash + "WMI Queries"), $true, $false, $false, $false,
    $false, $true, 1, $false, ("Task 3" + $emDash + "WMI Queries with CIM"), 2) | Out-Null
